$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to remain text even when the value looks numeric
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

function Set-DValue($ref, $value) {
    if ($value -match "^-?[0-9]+(\.[0-9]+)?$") {
        Set-TextValue $ws.Range($ref) $value
    } else {
        $ws.Range($ref).Value = $value
    }
}

# Rows 23 and 24 swap places (Uniswap moves above InternetComputer(DFINITY))
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-DValue "D23" "7.55"
$ws.Range("E23").Value = "  +11.00%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-DValue "D24" "14.26"
$ws.Range("E24").Value = "  +19.16%  "

# Remaining price / volume refreshes
Set-DValue "D2" "63.778.90"
$ws.Range("E2").Value = "  +6.01%  "
Set-DValue "D3" "3.054.71"
$ws.Range("E3").Value = "  +5.67%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-DValue "D5" "556.70"
$ws.Range("E5").Value = "  +5.31%  "
Set-DValue "D6" "142.53"
$ws.Range("E6").Value = "  +9.62%  "
$ws.Range("E7").Value = "  -0.14%  "
Set-DValue "D8" "3.052.26"
$ws.Range("E8").Value = "  +5.81%  "
$ws.Range("E9").Value = "  +7.42%  "
Set-DValue "D10" "0.156"
$ws.Range("E10").Value = "  +10.52%  "
Set-DValue "D11" "6.08"
$ws.Range("E11").Value = "  -3.82%  "
Set-DValue "D12" "0.479"
$ws.Range("E12").Value = "  +11.94%  "
$ws.Range("E13").Value = "  +9.51%  "
Set-DValue "D14" "35.02"
Set-DValue "D15" "3.550.69"
$ws.Range("E15").Value = "  +4.78%  "
Set-DValue "D16" "63.777.92"
$ws.Range("E16").Value = "  +6.11%  "
$ws.Range("E17").Value = "  +4.13%  "
Set-DValue "D18" "3.055.85"
$ws.Range("E18").Value = "  +5.70%  "
Set-DValue "D19" "6.74"
$ws.Range("E19").Value = "  +6.22%  "
Set-DValue "D20" "478.13"
$ws.Range("E20").Value = "  +6.41%  "
$ws.Range("E21").Value = "  +8.65%  "
Set-DValue "D22" "0.678"
$ws.Range("E22").Value = "  +8.21%  "
Set-DValue "D25" "81.40"
$ws.Range("E25").Value = "  +6.18%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +7.04%  "
Set-DValue "D28" "7.93"
$ws.Range("E28").Value = "  +9.38%  "
$ws.Range("E29").Value = "  +7.23%  "
Set-DValue "D30" "1.00"
$ws.Range("E30").Value = "  +0.19%  "
Set-DValue "D31" "26.21"
$ws.Range("E31").Value = "  +7.58%  "
$ws.Range("E32").Value = "  +4.02%  "
$ws.Range("E33").Value = "  +9.20%  "
Set-DValue "D34" "5.62"
$ws.Range("E34").Value = "  +6.36%  "
$ws.Range("E35").Value = "  +11.58%  "
Set-DValue "D36" "54.96"
$ws.Range("E36").Value = "  +3.43%  "
Set-DValue "D37" "0.0407"
$ws.Range("E37").Value = "  +9.69%  "
Set-DValue "D38" "445.10"
$ws.Range("E38").Value = "  +3.15%  "
$ws.Range("E39").Value = "  +5.11%  "
Set-DValue "D40" "2.83"
$ws.Range("E40").Value = "  +23.72%  "
Set-DValue "D41" "2.963.53"
$ws.Range("E41").Value = "  +3.82%  "
Set-DValue "D42" "8.24"
$ws.Range("E42").Value = "  +7.09%  "
$ws.Range("E43").Value = "  +1.96%  "
Set-DValue "D44" "27.74"
$ws.Range("E44").Value = "  +9.54%  "
$ws.Range("E45").Value = "  +10.37%  "
Set-DValue "D46" "2.16"
$ws.Range("E46").Value = "  +14.71%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +7.71%  "
Set-DValue "D49" "0.0₃0514"
$ws.Range("E49").Value = "  +10.09%  "
Set-DValue "D50" "116.74"
$ws.Range("E50").Value = "  +4.33%  "
Set-DValue "D51" "2.08"
$ws.Range("E51").Value = "  +9.34%  "
